# Updated CHE model - 2025-08-31 22:13
#
# The source model renamed the "wind" label used for the onshore-wind
# build-rate lookup row to "windon" (distinguishing it from the generic
# "wind" label still used elsewhere in the workbook). This script mirrors
# that edit:
#   - Veda!F9 (the onshore-wind row key of the F3:J11 lookup table) changes
#     from "wind" to "windon" - which also shifts the VLOOKUP-based D7
#     result (recalculated automatically).
#   - Every "wind" tag in the historical_data_long sheet's technology
#     column (A) that refers to the wind-onshore series is relabeled to
#     "windon" to match.

$wb = $excel.ActiveWorkbook

# --- Veda sheet: build-rate lookup table key ---
$wsVeda = $wb.Worksheets.Item("Veda")
$wsVeda.Range("F9").Value = "windon"

# --- historical_data_long sheet: relabel all "wind" rows to "windon" ---
$wsHist = $wb.Worksheets.Item("historical_data_long")
$windRows = @(10,18,26,34,42,50,58,66,74,82,90,98,106,114,122,130,138,146,154,162,170,178,186,194,202,210,218,226,234,242,250,258,266,274,282,290,298,306,314,322,330,338,346,354,362,370,378,386,394,402,410,418,426,434,442,450,458,466,474,482,490,498,506,514,522,530,538,546,554,562,570,578)

foreach ($r in $windRows) {
    $wsHist.Cells.Item($r, 1).Value = "windon"
}
